$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'97.863.85"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.09%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.424.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +4.03%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.10%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'255.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.31%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'654.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +4.44%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.49"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +2.22%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.429"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +6.87%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'1.07"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +8.38%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.05%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'3.422.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +4.05%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +4.07%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'41.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.41%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +15.33%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0000258"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.29%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'97.524.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.19%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.060.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.90%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'8.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +34.15%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.433.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.19%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +13.04%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.499"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +45.50%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'10.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +14.27%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'3.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.72%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'505.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +3.51%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0000206"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.45%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'6.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +8.19%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'98.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +11.19%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'12.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +4.41%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'3.601.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +4.01%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +4.07%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.202"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +6.17%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'11.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +6.05%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.24%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.06%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +19.34%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'29.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +6.72%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +15.76%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'7.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +6.12%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'Fetch.AI"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'1.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +15.14%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'Kaspa"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'0.154"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.17%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'520.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +4.95%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'24.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.876"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +11.78%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'3.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.58%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.0419"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +23.92%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'5.56"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +15.16%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +4.48%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'8.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +12.33%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.02%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +13.22%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +5.18%  "
$ws.Range("E51").Style = "Normal"
